# Fruta / hortaliza, semanal
# Insert 3 new weekly price rows for "Nectarín" at Vega Modelo de Temuco
# (La Araucanía) just before the existing row 473, pushing the rest of
# the table down by three rows (old 473..503 -> new 476..506).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows starting at row 473 (inserting the same row
# three times pushes everything below down by one each time).
$ws.Rows.Item(473).Insert()
$ws.Rows.Item(473).Insert()
$ws.Rows.Item(473).Insert()

# Common (constant) columns shared by every data row in this sheet.
$mercadoId = 10
$mercado = "Vega Modelo de Temuco"
$region = "La Araucanía"
$codreg = 9
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria = "Nectarín"

$newRows = @(
    @(44610, "June Pearl",    "Primera", 120, 17000, 17000, 17000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 944, 18),
    @(44610, "Ruby Diamond",  "Primera", 100, 17000, 17000, 17000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 944, 18),
    @(44610, "Venus",         "Primera", 200, 17000, 17000, 17000, "`$/bandeja 18 kilos granel", "Región de O'Higgins", 944, 18)
)

$r = 473
foreach ($row in $newRows) {
    $fecha = $row[0]
    $variedad = $row[1]
    $calidad = $row[2]
    $volumen = $row[3]
    $precioMin = $row[4]
    $precioMax = $row[5]
    $precioProm = $row[6]
    $unidad = $row[7]
    $origen = $row[8]
    $precioKg = $row[9]
    $kgUnidad = $row[10]

    $ws.Cells.Item($r, 1).Value = $mercadoId
    $ws.Cells.Item($r, 2).Value = $mercado
    $ws.Cells.Item($r, 3).Value = $region
    $ws.Cells.Item($r, 4).Value = $fecha
    $ws.Cells.Item($r, 5).Value = $codreg
    $ws.Cells.Item($r, 6).Value = $tipo
    $ws.Cells.Item($r, 7).Value = $productoId
    $ws.Cells.Item($r, 8).Value = $producto
    $ws.Cells.Item($r, 9).Value = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $variedad
    $ws.Cells.Item($r, 12).Value = $calidad
    $ws.Cells.Item($r, 13).Value = $volumen
    $ws.Cells.Item($r, 14).Value = $precioMin
    $ws.Cells.Item($r, 15).Value = $precioMax
    $ws.Cells.Item($r, 16).Value = $precioProm
    $ws.Cells.Item($r, 17).Value = $unidad
    $ws.Cells.Item($r, 18).Value = $origen
    $ws.Cells.Item($r, 19).Value = $precioKg
    $ws.Cells.Item($r, 20).Value = $kgUnidad

    $r = $r + 1
}
